$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("TextBox 9")

# Update position/size (EMU -> points: 914400 EMU = 72 points)
$shape.Left = 274.0656
$shape.Top = 488.0816
$shape.Width = 226.8709
$shape.Height = 50.89221

# Update text: add a new paragraph "(Named Kasa-obake)" after "Example of a Yokai"
$tf = $shape.TextFrame
$tr = $tf.TextRange
$tr.Text = "Example of a Yokai`r(Named Kasa-obake)"

# Re-assign the "obake" substring's text in place so it becomes its own run
# (mirrors the misspelling-flagged run in the authored deck).
$sub = $tr.Characters(32, 5)
$sub.Text = "obake"
